$d = $word.ActiveDocument

# --- Typo fix: "of earth, how fast will the station need to be spinning."
#     -> "of Earth, how fast will the station need to be spinning?"
#     (capitalize "Earth" and turn the trailing period into a question mark)

$find = $d.Content.Find
$find.ClearFormatting()
$find.Replacement.ClearFormatting()
$found = $find.Execute(
    "of earth, how fast will the station need to be spinning.",
    $true,
    $false,
    $false,
    $false,
    $false,
    $true,
    1,
    $false,
    "of Earth, how fast will the station need to be spinning?",
    2)

if (-not $found) {
    throw "Could not find the sentence to fix about the space station / Earth gravity."
}
